# Daily refresh of the cryptocurrency price/volume table on Sheet1.
# Mirrors the nightly GitHub Actions job that re-scrapes coinranking.com
# and rewrites columns D (Price) and E (Volume(1h)) for each ranked coin.
# A handful of coins (Algorand/WEMIXTOKEN, NEARProtocol/Decentraland) also
# swapped rank positions since the previous run, so their Coin name (B),
# Link (C), Price (D) and Volume (E) cells move together as whole rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values such as "1.002" or "10.80" read as plain numbers through
# Excel's normal text->value coercion, which would silently drop the
# significant trailing zero / digits baked into the scraped string.
# Forcing the cell to Text format before the assignment keeps the
# original string (including trailing zeros) exactly as scraped.
function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $Text
}

$ws.Range('D2').Value = '24.553.22'
$ws.Range('E2').Value = '  -0.68%  '
$ws.Range('D3').Value = '1.691.49'
$ws.Range('E3').Value = '  -0.18%  '
Set-TextValue 'D4' '1.002'
$ws.Range('E4').Value = '  -0.07%  '
Set-TextValue 'D5' '314.19'
$ws.Range('E5').Value = '  -0.88%  '
Set-TextValue 'D6' '1.002'
$ws.Range('E6').Value = '  -0.03%  '
Set-TextValue 'D7' '0.3890'
$ws.Range('E7').Value = '  -1.67%  '
Set-TextValue 'D8' '0.4036'
$ws.Range('E8').Value = '  -0.78%  '
$ws.Range('E9').Value = '  -0.08%  '
Set-TextValue 'D10' '1.006'
$ws.Range('E10').Value = '  +0.23%  '
Set-TextValue 'D11' '52.67'
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('E12').Value = '  -1.97%  '
Set-TextValue 'D13' '25.22'
$ws.Range('E13').Value = '  +6.80%  '
Set-TextValue 'D14' '7.516'
$ws.Range('E14').Value = '  +3.44%  '
Set-TextValue 'D15' '0.00001355'
$ws.Range('E15').Value = '  +2.40%  '
Set-TextValue 'D16' '7.970'
$ws.Range('E16').Value = '  -1.01%  '
$ws.Range('D17').Value = '1.682.64'
$ws.Range('E17').Value = '  -0.88%  '
Set-TextValue 'D18' '98.49'
$ws.Range('E18').Value = '  -1.50%  '
Set-TextValue 'D19' '0.07104'
$ws.Range('E19').Value = '  +0.90%  '
Set-TextValue 'D20' '20.02'
$ws.Range('E20').Value = '  +1.88%  '
Set-TextValue 'D21' '7.285'
$ws.Range('E21').Value = '  +4.02%  '
Set-TextValue 'D22' '1.003'
Set-TextValue 'D23' '14.26'
$ws.Range('E23').Value = '  -0.57%  '
$ws.Range('D24').Value = '24.542.60'
$ws.Range('E24').Value = '  -0.62%  '
Set-TextValue 'D25' '2.979'
$ws.Range('E25').Value = '  -9.33%  '
Set-TextValue 'D26' '2.351'
$ws.Range('E26').Value = '  -0.35%  '
$ws.Range('E27').Value = '  -0.06%  '
Set-TextValue 'D28' '162.18'
$ws.Range('E28').Value = '  -0.08%  '
Set-TextValue 'D29' '8.741'
$ws.Range('E29').Value = '  +16.26%  '
Set-TextValue 'D30' '136.94'
$ws.Range('E30').Value = '  +0.54%  '
Set-TextValue 'D31' '5.223'
$ws.Range('E31').Value = '  +0.45%  '
$ws.Range('D32').Value = '1.869.88'
$ws.Range('E32').Value = '  -0.73%  '
Set-TextValue 'D33' '0.08845'
$ws.Range('E33').Value = '  +2.04%  '
Set-TextValue 'D34' '7.406'
$ws.Range('E34').Value = '  +4.45%  '
Set-TextValue 'D35' '1.033'
$ws.Range('E35').Value = '  -2.12%  '
$ws.Range('B36').Value = 'Algorand'
$ws.Range('C36').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D36' '0.2792'
$ws.Range('E36').Value = '  +1.97%  '
$ws.Range('B37').Value = 'WEMIXTOKEN'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D37' '1.966'
$ws.Range('E37').Value = '  +4.37%  '
Set-TextValue 'D38' '0.02918'
$ws.Range('E38').Value = '  +6.95%  '
Set-TextValue 'D39' '10.80'
$ws.Range('E39').Value = '  -5.41%  '
Set-TextValue 'D40' '14.25'
$ws.Range('E40').Value = '  -1.62%  '
Set-TextValue 'D41' '0.09142'
$ws.Range('E41').Value = '  -1.19%  '
Set-TextValue 'D42' '0.7945'
$ws.Range('E42').Value = '  +3.51%  '
Set-TextValue 'D43' '1.460'
$ws.Range('E43').Value = '  -1.01%  '
Set-TextValue 'D44' '16.71'
$ws.Range('E44').Value = '  +3.36%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D45' '2.615'
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D46' '0.7236'
$ws.Range('E46').Value = '  +0.81%  '
Set-TextValue 'D47' '4.198'
$ws.Range('E47').Value = '  -0.50%  '
Set-TextValue 'D49' '1.351'
$ws.Range('E49').Value = '  +2.26%  '
Set-TextValue 'D50' '138.19'
$ws.Range('E50').Value = '  -1.61%  '
Set-TextValue 'D51' '91.14'
$ws.Range('E51').Value = '  +0.25%  '
